# Update BunkerPrices at 2025-03-28 02:23
#
# 1) Reorder the "Khor Fakkan" column: it moves from column AU to column X
#    (inserted just before "Petersburg"); the columns Petersburg, Sydney,
#    Santos, Istanbul and "South Korea (West)" each shift one place to the
#    right (X->Y->Z->AA->AB->AC), and the column that used to be at AC
#    ("Antwerp") ends up at AU, where "Khor Fakkan" used to live. All other
#    columns (AD..AT and AV) are untouched.
# 2) The "Date" column's number format on the (previously) last data row
#    (row 10) reverts to the normal "YYYY-MM-DD HH:MM:SS" format used by
#    every other data row.
# 3) A brand-new data row (row 11, dated 2025-03-25) is appended, with its
#    "Date" cell using the special "YYYY-MM-DD" format reserved for the
#    most-recent row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate the seven affected columns across header + all data rows ---
$rotatedCols = @(24, 25, 26, 27, 28, 29, 47)   # X, Y, Z, AA, AB, AC, AU

for ($r = 1; $r -le 10; $r++) {
    $oldVals = @{}
    foreach ($c in $rotatedCols) {
        $oldVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $ws.Cells.Item($r, 24).Value = $oldVals[47]   # X  <- old AU (Khor Fakkan)
    $ws.Cells.Item($r, 25).Value = $oldVals[24]   # Y  <- old X  (Petersburg)
    $ws.Cells.Item($r, 26).Value = $oldVals[25]   # Z  <- old Y  (Sydney)
    $ws.Cells.Item($r, 27).Value = $oldVals[26]   # AA <- old Z  (Santos)
    $ws.Cells.Item($r, 28).Value = $oldVals[27]   # AB <- old AA (Istanbul)
    $ws.Cells.Item($r, 29).Value = $oldVals[28]   # AC <- old AB (South Korea (West))
    $ws.Cells.Item($r, 47).Value = $oldVals[29]   # AU <- old AC (Antwerp)
}

# --- Step 2: row 10's Date cell goes back to the regular datetime format ---
$ws.Cells.Item(10, 19).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Step 3: append the new row 11 ---
$row11 = @{
    1=520; 2=572; 3=535; 4=517; 5=650; 6=522; 7=650; 8=659; 9=582; 10=520;
    11=562; 12=513; 13=585; 14=523; 15=645; 16=767; 17=570; 18=665;
    19=45741;
    20=585; 21=592; 22=607; 23=508; 24=513; 25=555; 26=757; 27=538; 28=574;
    29=527; 30=651; 31=603.5; 32=562; 33=528; 34=573; 35=882; 36=650; 37=507;
    38=623; 39=547; 40=517; 41=535; 42=513; 43=508; 44=490; 45=522; 46=552;
    47=490; 48=555
}

foreach ($c in $row11.Keys) {
    $ws.Cells.Item(11, $c).Value = $row11[$c]
}

# The newest row's Date cell uses the special short date-only format.
$ws.Cells.Item(11, 19).NumberFormat = "YYYY-MM-DD"
